$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rng = $ws.Range("B1:B1048576")

$rng.FormatConditions.Delete()

$fc1 = $rng.FormatConditions.Add(1, 3, '"no comenzado"')
$fc1.Interior.Color = 255
$fc2 = $rng.FormatConditions.Add(1, 3, '"en proceso"')
$fc2.Interior.Color = 65535
$fc3 = $rng.FormatConditions.Add(1, 3, '"terminado"')
$fc3.Interior.Color = 5296274

for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $fc = $rng.FormatConditions.Item($i)
    Write-Host "FC $i : $($fc.Formula1) color=$($fc.Interior.Color)"
}
